$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 40.75339133333333
$ws.Range("H2").Value = 122.260174
$ws.Range("I2").Value = 0.02126536631186857
$ws.Range("J2").Value = 0.02126536631186857
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.8426706666666667
$ws.Range("N2").Value = 2.528012
$ws.Range("O2").Value = 0.07312747333774275
$ws.Range("P2").Value = 0.07312747333774275
$ws.Range("Q2").Value = 34.34168744378756
$ws.Range("R2").Value = 309.075186994088
$ws.Range("S2").Value = 0.001555082507988501
$ws.Range("T2").Value = 0.001555082507988501
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 40.75339133333333
$ws.Range("H3").Value = 122.260174
$ws.Range("I3").Value = 0.02126536631186857
$ws.Range("J3").Value = 0.02126536631186857
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.843693666666667
$ws.Range("N3").Value = 8.531081
$ws.Range("O3").Value = 0.2467774671835513
$ws.Range("P3").Value = 0.2467774671835513
$ws.Range("Q3").Value = 115.8901608297882
$ws.Range("R3").Value = 1043.011447468094
$ws.Range("S3").Value = 0.005247813237173342
$ws.Range("T3").Value = 0.005247813237173341
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 40.75339133333333
$ws.Range("H4").Value = 122.260174
$ws.Range("I4").Value = 0.02126536631186857
$ws.Range("J4").Value = 0.02126536631186857
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.430413
$ws.Range("N4").Value = 1.291239
$ws.Range("O4").Value = 0.03735150210725013
$ws.Range("P4").Value = 0.03735150210725013
$ws.Range("Q4").Value = 17.540789423954
$ws.Range("R4").Value = 157.867104815586
$ws.Range("S4").Value = 0.0007942933746092046
$ws.Range("T4").Value = 0.0007942933746092046
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 40.75339133333333
$ws.Range("H5").Value = 122.260174
$ws.Range("I5").Value = 0.02126536631186857
$ws.Range("J5").Value = 0.02126536631186857
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 7.406534333333333
$ws.Range("N5").Value = 22.219603
$ws.Range("O5").Value = 0.6427435573714558
$ws.Range("P5").Value = 0.6427435573714558
$ws.Range("Q5").Value = 301.8413921101024
$ws.Range("R5").Value = 2716.572528990922
$ws.Range("S5").Value = 0.01366817719209752
$ws.Range("T5").Value = 0.01366817719209752
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1689.289306666667
$ws.Range("H6").Value = 5067.86792
$ws.Range("I6").Value = 0.8814813868902838
$ws.Range("J6").Value = 0.8814813868902838
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.8426706666666667
$ws.Range("N6").Value = 2.528012
$ws.Range("O6").Value = 0.07312747333774275
$ws.Range("P6").Value = 0.07312747333774275
$ws.Range("Q6").Value = 1423.514546241671
$ws.Range("R6").Value = 12811.63091617504
$ws.Range("S6").Value = 0.06446050661753573
$ws.Range("T6").Value = 0.06446050661753573
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1689.289306666667
$ws.Range("H7").Value = 5067.86792
$ws.Range("I7").Value = 0.8814813868902838
$ws.Range("J7").Value = 0.8814813868902838
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.843693666666667
$ws.Range("N7").Value = 8.531081
$ws.Range("O7").Value = 0.2467774671835513
$ws.Range("P7").Value = 0.2467774671835513
$ws.Range("Q7").Value = 4803.821302535724
$ws.Range("R7").Value = 43234.39172282152
$ws.Range("S7").Value = 0.2175297440262283
$ws.Range("T7").Value = 0.2175297440262283
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1689.289306666667
$ws.Range("H8").Value = 5067.86792
$ws.Range("I8").Value = 0.8814813868902838
$ws.Range("J8").Value = 0.8814813868902838
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.430413
$ws.Range("N8").Value = 1.291239
$ws.Range("O8").Value = 0.03735150210725013
$ws.Range("P8").Value = 0.03735150210725013
$ws.Range("Q8").Value = 727.09207835032
$ws.Range("R8").Value = 6543.82870515288
$ws.Range("S8").Value = 0.0329246538799342
$ws.Range("T8").Value = 0.0329246538799342
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1689.289306666667
$ws.Range("H9").Value = 5067.86792
$ws.Range("I9").Value = 0.8814813868902838
$ws.Range("J9").Value = 0.8814813868902838
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 7.406534333333333
$ws.Range("N9").Value = 22.219603
$ws.Range("O9").Value = 0.6427435573714558
$ws.Range("P9").Value = 0.6427435573714558
$ws.Range("Q9").Value = 12511.77924875953
$ws.Range("R9").Value = 112606.0132388357
$ws.Range("S9").Value = 0.5665664823665856
$ws.Range("T9").Value = 0.5665664823665856
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 100.9654023333333
$ws.Range("H10").Value = 302.896207
$ws.Range("I10").Value = 0.05268435816499466
$ws.Range("J10").Value = 0.05268435816499466
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.8426706666666667
$ws.Range("N10").Value = 2.528012
$ws.Range("O10").Value = 0.07312747333774275
$ws.Range("P10").Value = 0.07312747333774275
$ws.Range("Q10").Value = 85.08058289449822
$ws.Range("R10").Value = 765.725246050484
$ws.Range("S10").Value = 0.003852673997026736
$ws.Range("T10").Value = 0.003852673997026736
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 100.9654023333333
$ws.Range("H11").Value = 302.896207
$ws.Range("I11").Value = 0.05268435816499466
$ws.Range("J11").Value = 0.05268435816499466
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 2.843693666666667
$ws.Range("N11").Value = 8.531081
$ws.Range("O11").Value = 0.2467774671835513
$ws.Range("P11").Value = 0.2467774671835513
$ws.Range("Q11").Value = 287.1146751677519
$ws.Range("R11").Value = 2584.032076509767
$ws.Range("S11").Value = 0.01300131246814843
$ws.Range("T11").Value = 0.01300131246814843
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 100.9654023333333
$ws.Range("H12").Value = 302.896207
$ws.Range("I12").Value = 0.05268435816499466
$ws.Range("J12").Value = 0.05268435816499466
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.430413
$ws.Range("N12").Value = 1.291239
$ws.Range("O12").Value = 0.03735150210725013
$ws.Range("P12").Value = 0.03735150210725013
$ws.Range("Q12").Value = 43.456821714497
$ws.Range("R12").Value = 391.111395430473
$ws.Range("S12").Value = 0.001967839915018919
$ws.Range("T12").Value = 0.001967839915018919
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 100.9654023333333
$ws.Range("H13").Value = 302.896207
$ws.Range("I13").Value = 0.05268435816499466
$ws.Range("J13").Value = 0.05268435816499466
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 7.406534333333333
$ws.Range("N13").Value = 22.219603
$ws.Range("O13").Value = 0.6427435573714558
$ws.Range("P13").Value = 0.6427435573714558
$ws.Range("Q13").Value = 747.8037188606468
$ws.Range("R13").Value = 6730.233469745821
$ws.Range("S13").Value = 0.03386253178480057
$ws.Range("T13").Value = 0.03386253178480057
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 85.41274733333334
$ws.Range("H14").Value = 256.238242
$ws.Range("I14").Value = 0.04456888863285297
$ws.Range("J14").Value = 0.04456888863285297
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 0.8426706666666667
$ws.Range("N14").Value = 2.528012
$ws.Range("O14").Value = 0.07312747333774275
$ws.Range("P14").Value = 0.07312747333774275
$ws.Range("Q14").Value = 71.97481673721157
$ws.Range("R14").Value = 647.773350634904
$ws.Range("S14").Value = 0.003259210215191781
$ws.Range("T14").Value = 0.003259210215191781
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 85.41274733333334
$ws.Range("H15").Value = 256.238242
$ws.Range("I15").Value = 0.04456888863285297
$ws.Range("J15").Value = 0.04456888863285297
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 2.843693666666667
$ws.Range("N15").Value = 8.531081
$ws.Range("O15").Value = 0.2467774671835513
$ws.Range("P15").Value = 0.2467774671835513
$ws.Range("Q15").Value = 242.8876886444003
$ws.Range("R15").Value = 2185.989197799602
$ws.Range("S15").Value = 0.01099859745200123
$ws.Range("T15").Value = 0.01099859745200122
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 85.41274733333334
$ws.Range("H16").Value = 256.238242
$ws.Range("I16").Value = 0.04456888863285297
$ws.Range("J16").Value = 0.04456888863285297
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.430413
$ws.Range("N16").Value = 1.291239
$ws.Range("O16").Value = 0.03735150210725013
$ws.Range("P16").Value = 0.03735150210725013
$ws.Range("Q16").Value = 36.762756817982
$ws.Range("R16").Value = 330.864811361838
$ws.Range("S16").Value = 0.001664714937687804
$ws.Range("T16").Value = 0.001664714937687804
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 85.41274733333334
$ws.Range("H17").Value = 256.238242
$ws.Range("I17").Value = 0.04456888863285297
$ws.Range("J17").Value = 0.04456888863285297
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 7.406534333333333
$ws.Range("N17").Value = 22.219603
$ws.Range("O17").Value = 0.6427435573714558
$ws.Range("P17").Value = 0.6427435573714558
$ws.Range("Q17").Value = 632.6124456286585
$ws.Range("R17").Value = 5693.512010657926
$ws.Range("S17").Value = 0.02864636602797216
$ws.Range("T17").Value = 0.02864636602797216